$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) from 45177 (2023-09-08) to
# 45178 (2023-09-09) for every data row (rows 2 through 300).
$ws.Range("C2:C300").Value = 45178
